$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 850
$ws.Range("I9").Value = 850
$ws.Range("K9").Value = 850
$ws.Range("M9").Value = -681
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1500
$ws.Range("K29").Value = 4500
$ws.Range("M29").Value = -4219
$ws.Range("H98").Value = 5453.826
$ws.Range("I98").Value = 5247.636
$ws.Range("K98").Value = 5247.636
$ws.Range("M98").Value = -3749.636
$ws.Range("H99").Value = 500
$ws.Range("I99").Value = 800
$ws.Range("K99").Value = 2400
$ws.Range("M99").Value = -902
$ws.Range("H100").Value = 1852.3889
$ws.Range("I100").Value = 1299.1538
$ws.Range("J100").Value = 3290.8
$ws.Range("K100").Value = 1299.1538
$ws.Range("L100").Value = 3290.8
$ws.Range("M100").Value = -758.1538
$ws.Range("N100").Value = -4372.8
$ws.Range("H103").Value = 619.86365
$ws.Range("J103").Value = 1236.625
$ws.Range("L103").Value = 3709.875
$ws.Range("N103").Value = -4881.875
$ws.Range("H122").Value = 5453.826
$ws.Range("I122").Value = 5247.636
$ws.Range("K122").Value = 15742.908
$ws.Range("M122").Value = -13292.908
$ws.Range("H137").Value = 1808
$ws.Range("I137").Value = 1918.6842
$ws.Range("J137").Value = 1632.75
$ws.Range("K137").Value = 5756.0526
$ws.Range("L137").Value = 4898.25
$ws.Range("M137").Value = -3206.0526
$ws.Range("N137").Value = -9998.25
$ws.Range("H138").Value = 6147.577
$ws.Range("J138").Value = 15833.223
$ws.Range("L138").Value = 47499.669
$ws.Range("N138").Value = -57779.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4819.068
$ws.Range("I61").Value = 2199.4666
$ws.Range("J61").Value = 13239.214
$ws.Range("K61").Value = 2199.4666
$ws.Range("L61").Value = 13239.214
$ws.Range("M61").Value = -1987.4666
$ws.Range("N61").Value = -13663.214
$ws.Range("H74").Value = 71439.44500000001
$ws.Range("I74").Value = 301990
$ws.Range("J74").Value = 5567.857
$ws.Range("K74").Value = 301990
$ws.Range("L74").Value = 5567.857
$ws.Range("M74").Value = -301116
$ws.Range("N74").Value = -7315.857
$ws.Range("H77").Value = 71439.44500000001
$ws.Range("I77").Value = 301990
$ws.Range("J77").Value = 5567.857
$ws.Range("K77").Value = 1509950
$ws.Range("L77").Value = 27839.285
$ws.Range("M77").Value = -1505582
$ws.Range("N77").Value = -36575.285
$ws.Range("H82").Value = 47160.5
$ws.Range("J82").Value = 47160.5
$ws.Range("L82").Value = 47160.5
$ws.Range("N82").Value = -47882.5
$ws.Range("H85").Value = 47160.5
$ws.Range("J85").Value = 47160.5
$ws.Range("L85").Value = 47160.5
$ws.Range("N85").Value = -49656.5
$ws.Range("H97").Value = 4636130.5
$ws.Range("I97").Value = 366.66666
$ws.Range("J97").Value = 27814948
$ws.Range("K97").Value = 366.66666
$ws.Range("L97").Value = 27814948
$ws.Range("M97").Value = 129.33334
$ws.Range("N97").Value = -27815940
$ws.Range("H136").Value = 4819.068
$ws.Range("I136").Value = 2199.4666
$ws.Range("J136").Value = 13239.214
$ws.Range("K136").Value = 6598.399800000001
$ws.Range("L136").Value = 39717.642
$ws.Range("M136").Value = -4048.399800000001
$ws.Range("N136").Value = -44817.642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4527.385
$ws.Range("I94").Value = 1104.5
$ws.Range("J94").Value = 10004
$ws.Range("K94").Value = 1104.5
$ws.Range("L94").Value = 10004
$ws.Range("M94").Value = -653.5
$ws.Range("N94").Value = -10906
$ws.Range("H107").Value = 41671736
$ws.Range("I107").Value = 48917388
$ws.Range("J107").Value = 9244.25
$ws.Range("K107").Value = 48917388
$ws.Range("L107").Value = 9244.25
$ws.Range("M107").Value = -48915468
$ws.Range("N107").Value = -13084.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11236.4
$ws.Range("I31").Value = 6572.2
$ws.Range("J31").Value = 12169.24
$ws.Range("K31").Value = 6572.2
$ws.Range("L31").Value = 12169.24
$ws.Range("M31").Value = -6277.2
$ws.Range("N31").Value = -12759.24
$ws.Range("H34").Value = 11236.4
$ws.Range("I34").Value = 6572.2
$ws.Range("J34").Value = 12169.24
$ws.Range("K34").Value = 6572.2
$ws.Range("L34").Value = 12169.24
$ws.Range("M34").Value = -6370.2
$ws.Range("N34").Value = -12573.24
$ws.Range("H58").Value = 13895250
$ws.Range("I58").Value = 35715070
$ws.Range("J58").Value = 9906.682000000001
$ws.Range("K58").Value = 35715070
$ws.Range("L58").Value = 9906.682000000001
$ws.Range("M58").Value = -35714867
$ws.Range("N58").Value = -10312.682
$ws.Range("H136").Value = 13895250
$ws.Range("I136").Value = 35715070
$ws.Range("J136").Value = 9906.682000000001
$ws.Range("K136").Value = 107145210
$ws.Range("L136").Value = 29720.046
$ws.Range("M136").Value = -107142660
$ws.Range("N136").Value = -34820.046

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 707.3570999999999
$ws.Range("J12").Value = 992
$ws.Range("L12").Value = 2976
$ws.Range("N12").Value = -3322
$ws.Range("H75").Value = 55559624
$ws.Range("I75").Value = 83336060
$ws.Range("K75").Value = 250008180
$ws.Range("M75").Value = -250007182
$ws.Range("H78").Value = 55559624
$ws.Range("I78").Value = 83336060
$ws.Range("K78").Value = 750024540
$ws.Range("M78").Value = -750019548
$ws.Range("H98").Value = 1167.7142
$ws.Range("J98").Value = 1427.2858
$ws.Range("L98").Value = 4281.857400000001
$ws.Range("N98").Value = -7277.857400000001
$ws.Range("H103").Value = 307.2
$ws.Range("I103").Value = 230
$ws.Range("J103").Value = 487.33334
$ws.Range("K103").Value = 690
$ws.Range("L103").Value = 1462.00002
$ws.Range("M103").Value = 189
$ws.Range("N103").Value = -3220.00002
$ws.Range("H132").Value = 9851.191999999999
$ws.Range("I132").Value = 4224.857
$ws.Range("J132").Value = 16415.25
$ws.Range("K132").Value = 38023.713
$ws.Range("L132").Value = 147737.25
$ws.Range("M132").Value = -35493.713
$ws.Range("N132").Value = -152797.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9080.091
$ws.Range("I102").Value = 10411.714
$ws.Range("J102").Value = 6749.75
$ws.Range("K102").Value = 10411.714
$ws.Range("L102").Value = 6749.75
$ws.Range("M102").Value = -8789.714
$ws.Range("N102").Value = -9993.75
$ws.Range("H132").Value = 11208.111
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 13196.143
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 39588.429
$ws.Range("M132").Value = -10220
$ws.Range("N132").Value = -44648.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4842.143
$ws.Range("I16").Value = 4173.875
$ws.Range("K16").Value = 4173.875
$ws.Range("M16").Value = -4003.875
$ws.Range("H61").Value = 3451714
$ws.Range("I61").Value = 5883436
$ws.Range("K61").Value = 5883436
$ws.Range("M61").Value = -5883234
$ws.Range("H93").Value = 2379.75
$ws.Range("I93").Value = 2271.8333
$ws.Range("J93").Value = 2703.5
$ws.Range("K93").Value = 2271.8333
$ws.Range("L93").Value = 2703.5
$ws.Range("M93").Value = -1023.8333
$ws.Range("N93").Value = -5199.5
$ws.Range("H113").Value = 3451714
$ws.Range("I113").Value = 5883436
$ws.Range("K113").Value = 5883436
$ws.Range("M113").Value = -5881266
$ws.Range("H122").Value = 5876.952
$ws.Range("I122").Value = 5047
$ws.Range("K122").Value = 15141
$ws.Range("M122").Value = -12691
$ws.Range("H132").Value = 10006106
$ws.Range("I132").Value = 21741730
$ws.Range("J132").Value = 9092.888999999999
$ws.Range("K132").Value = 65225190
$ws.Range("L132").Value = 27278.667
$ws.Range("M132").Value = -65222660
$ws.Range("N132").Value = -32338.667
